# Update generated "想去人数" (attendance) figures on the 展览 and 全部类型 sheets
# to reflect refreshed output, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1226
$ws1.Range("F5").Value  = 891
$ws1.Range("F6").Value  = 1626
$ws1.Range("F8").Value  = 1099
$ws1.Range("F19").Value = 35
$ws1.Range("F20").Value = 612
$ws1.Range("F21").Value = 603
$ws1.Range("F25").Value = 276
$ws1.Range("F26").Value = 5
$ws1.Range("F27").Value = 219

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1226
$ws4.Range("F6").Value  = 891
$ws4.Range("F7").Value  = 1626
$ws4.Range("F9").Value  = 1099
$ws4.Range("F25").Value = 35
$ws4.Range("F26").Value = 612
$ws4.Range("F27").Value = 603
$ws4.Range("F31").Value = 276
$ws4.Range("F33").Value = 5
$ws4.Range("F34").Value = 219

$wb.Save()
